$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Falta costo de CCC en documento estimación (now closed, no changes stage)
$ws.Range("E4").Value = 42331
$ws.Range("E4").NumberFormat = "DD/MM/YY"
$ws.Range("F4").Value = "Cerrada"
$ws.Range("G4").Value = "Se omite debido a que no se requiere la etapa de cambios para el proyecto"

# Row 5: Falta agregar la complejidad en tiempo de implementación (now closed)
$ws.Range("E5").Value = 42331
$ws.Range("E5").NumberFormat = "DD/MM/YY"
$ws.Range("F5").Value = "Cerrada"

# Row 6: No se cuenta con evidencia de envío de carta de aceptación al cliente (now closed)
$ws.Range("E6").Value = 42331
$ws.Range("E6").NumberFormat = "DD/MM/YY"
$ws.Range("F6").Value = "Cerrada"
$ws.Range("G6").Value = "Se muestra evidencia en Deal del cliente"

# Row 7: renamed non-conformity topic
$ws.Range("B7").Value = "Notificación de creación lineas base"

# Row 8: Falta evidencia de fechas reales de hitos y entregables (now closed)
$ws.Range("E8").Value = 42331
$ws.Range("E8").NumberFormat = "DD/MM/YY"
$ws.Range("F8").Value = "Cerrada"

# Row 4's comment is now longer and wraps onto more lines
$ws.Rows.Item(4).RowHeight = 41.75

# Leave the cursor where the author last left it
$ws.Range("G8").Select() | Out-Null
